$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.785.18'
$ws.Range("E2").Value = '  -2.41%  '

$ws.Range("D3").Value = '3.765.14'
$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '594.82'
$ws.Range("E5").Value = '  -3.06%  '

$ws.Range("D6").Value = '172.36'
$ws.Range("E6").Value = '  -2.93%  '

$ws.Range("D7").Value = '3.774.84'
$ws.Range("E7").Value = '  +0.69%  '

$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = '0.518'
$ws.Range("E9").Value = '  -1.60%  '

$ws.Range("D10").Value = '0.158'
$ws.Range("E10").Value = '  -5.29%  '

$ws.Range("D11").Value = '6.18'
$ws.Range("E11").Value = '  -5.72%  '

$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  -4.42%  '

$ws.Range("D13").Value = '37.59'
$ws.Range("E13").Value = '  -5.57%  '

$ws.Range("D14").Value = '0.0000242'
$ws.Range("E14").Value = '  -4.39%  '

$ws.Range("D15").Value = '4.382.81'
$ws.Range("E15").Value = '  +0.27%  '

$ws.Range("D16").Value = '3.751.71'
$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '67.786.89'
$ws.Range("E17").Value = '  -2.54%  '

$ws.Range("D18").Value = '0.115'
$ws.Range("E18").Value = '  -4.79%  '

$ws.Range("D19").Value = '7.16'
$ws.Range("E19").Value = '  -3.70%  '

$ws.Range("E20").Value = '  -1.41%  '

$ws.Range("D21").Value = '488.89'
$ws.Range("E21").Value = '  -2.14%  '

$ws.Range("D22").Value = '9.05'
$ws.Range("E22").Value = '  -1.25%  '

$ws.Range("D23").Value = '0.718'
$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").Value = '83.62'
$ws.Range("E24").Value = '  -2.47%  '

$ws.Range("D25").Value = '2.34'
$ws.Range("E25").Value = '  -8.89%  '

$ws.Range("D26").Value = '0.0000143'
$ws.Range("E26").Value = '  +6.99%  '

$ws.Range("D27").Value = '12.14'
$ws.Range("E27").Value = '  -5.80%  '

$ws.Range("E28").Value = '  -7.79%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").Value = '2.95'
$ws.Range("E30").Value = '  +1.49%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '32.68'
$ws.Range("E31").Value = '  +7.63%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '2.39'
$ws.Range("E32").Value = '  -3.61%  '

$ws.Range("D33").Value = '7.64'
$ws.Range("E33").Value = '  -5.09%  '

$ws.Range("E34").Value = '  -4.26%  '

$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").Value = '1.02'
$ws.Range("E36").Value = '  -2.58%  '

$ws.Range("E37").Value = '  -2.10%  '

$ws.Range("D38").Value = '5.71'
$ws.Range("E38").Value = '  -6.39%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '453.65'
$ws.Range("E39").Value = '  +2.53%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = '0.324'
$ws.Range("E40").Value = '  -6.71%  '

$ws.Range("D41").Value = '49.06'
$ws.Range("E41").Value = '  -1.26%  '

$ws.Range("E42").Value = '  -3.54%  '

$ws.Range("E43").Value = '  -7.32%  '

$ws.Range("D44").Value = '8.26'
$ws.Range("E44").Value = '  -3.26%  '

$ws.Range("D45").Value = '41.13'
$ws.Range("E45").Value = '  -7.51%  '

$ws.Range("D46").Value = '140.65'
$ws.Range("E46").Value = '  +1.68%  '

$ws.Range("D47").Value = '2.795.53'
$ws.Range("E47").Value = '  -5.21%  '

$ws.Range("D49").Value = '0.0347'
$ws.Range("E49").Value = '  -3.03%  '

$ws.Range("D50").Value = '25.69'
$ws.Range("E50").Value = '  -4.95%  '

$ws.Range("D51").Value = '23.31'
$ws.Range("E51").Value = '  +9.44%  '
